$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Left table (A:H), rows 1-6 only now (rows 7-8 cleared) ---
$ws.Range("A1").Value = "negative"

$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"

$ws.Range("A3").Value = "crude"
$ws.Range("B3").Value = 0.9117647058823529
$ws.Range("C3").Value = 31
$ws.Range("D3").Value = 31
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 3

$ws.Range("A4").Value = "crisis"
$ws.Range("B4").Value = 0.589041095890411
$ws.Range("C4").Value = 172
$ws.Range("D4").Value = 172
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 120

$ws.Range("A5").Value = "panic"
$ws.Range("B5").Value = 0.1724806201550388
$ws.Range("C5").Value = 89
$ws.Range("D5").Value = 89
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 427

$ws.Range("A6").Value = "sc"
$ws.Range("B6").Value = 0.1322751322751323
$ws.Range("C6").Value = 25
$ws.Range("D6").Value = 25
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 164

# Rows 7-8 of the left table no longer exist; clear them
$ws.Range("A7:H8").Clear()

# --- Right table (J:Q), rows 1-28 ---
$ws.Range("J1").Value = "positive"

$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"

$ws.Range("J3").Value = "interesting"
$ws.Range("K3").Value = 0.9393939393939394
$ws.Range("L3").Value = 31
$ws.Range("M3").Value = 31
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 2

$ws.Range("J4").Value = "love"
$ws.Range("K4").Value = 0.9347826086956522
$ws.Range("L4").Value = 43
$ws.Range("M4").Value = 43
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 3

$ws.Range("J5").Value = "best"
$ws.Range("K5").Value = 0.9152542372881356
$ws.Range("L5").Value = 54
$ws.Range("M5").Value = 54
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 5

$ws.Range("J6").Value = "great"
$ws.Range("K6").Value = 0.8660714285714286
$ws.Range("L6").Value = 97
$ws.Range("M6").Value = 97
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 15

$ws.Range("J7").Value = "thank"
$ws.Range("K7").Value = 0.8203125
$ws.Range("L7").Value = 105
$ws.Range("M7").Value = 105
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 23

$ws.Range("J8").Value = "thanks"
$ws.Range("K8").Value = 0.8048780487804879
$ws.Range("L8").Value = 66
$ws.Range("M8").Value = 66
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 16

$ws.Range("J9").Value = "free"
$ws.Range("K9").Value = 0.8
$ws.Range("L9").Value = 96
$ws.Range("M9").Value = 96
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 24

$ws.Range("J10").Value = "special"
$ws.Range("K10").Value = 0.7777777777777778
$ws.Range("L10").Value = 28
$ws.Range("M10").Value = 28
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 8

$ws.Range("J11").Value = "positive"
$ws.Range("K11").Value = 0.7758620689655172
$ws.Range("L11").Value = 45
$ws.Range("M11").Value = 45
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 13

$ws.Range("J12").Value = "safety"
$ws.Range("K12").Value = 0.7254901960784313
$ws.Range("L12").Value = 37
$ws.Range("M12").Value = 37
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 14

$ws.Range("J13").Value = "good"
$ws.Range("K13").Value = 0.71875
$ws.Range("L13").Value = 115
$ws.Range("M13").Value = 115
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 45

$ws.Range("J14").Value = "safe"
$ws.Range("K14").Value = 0.7112676056338029
$ws.Range("L14").Value = 101
$ws.Range("M14").Value = 101
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 41

$ws.Range("J15").Value = "confidence"
$ws.Range("K15").Value = 0.6944444444444444
$ws.Range("L15").Value = 25
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 11

$ws.Range("J16").Value = "support"
$ws.Range("K16").Value = 0.6886792452830188
$ws.Range("L16").Value = 73
$ws.Range("M16").Value = 73
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 33

$ws.Range("J17").Value = "relief"
$ws.Range("K17").Value = 0.64
$ws.Range("L17").Value = 32
$ws.Range("M17").Value = 32
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 18

$ws.Range("J18").Value = "well"
$ws.Range("K18").Value = 0.6276595744680851
$ws.Range("L18").Value = 59
$ws.Range("M18").Value = 59
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 35

$ws.Range("J19").Value = "better"
$ws.Range("K19").Value = 0.6190476190476191
$ws.Range("L19").Value = 39
$ws.Range("M19").Value = 39
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 24

$ws.Range("J20").Value = "fresh"
$ws.Range("K20").Value = 0.5833333333333334
$ws.Range("L20").Value = 28
$ws.Range("M20").Value = 28
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 20

$ws.Range("J21").Value = "heroes"
$ws.Range("K21").Value = 0.5531914893617021
$ws.Range("L21").Value = 26
$ws.Range("M21").Value = 26
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 21

$ws.Range("J22").Value = "hand"
$ws.Range("K22").Value = 0.5065274151436031
$ws.Range("L22").Value = 194
$ws.Range("M22").Value = 194
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 189

$ws.Range("J23").Value = "like"
$ws.Range("K23").Value = 0.4352941176470588
$ws.Range("L23").Value = 148
$ws.Range("M23").Value = 148
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 192

$ws.Range("J24").Value = "help"
$ws.Range("K24").Value = 0.4305084745762712
$ws.Range("L24").Value = 127
$ws.Range("M24").Value = 127
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 168

$ws.Range("J25").Value = "care"
$ws.Range("K25").Value = 0.4269662921348314
$ws.Range("L25").Value = 38
$ws.Range("M25").Value = 38
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 51

$ws.Range("J26").Value = "protect"
$ws.Range("K26").Value = 0.3972602739726027
$ws.Range("L26").Value = 29
$ws.Range("M26").Value = 29
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 44

$ws.Range("J27").Value = "increase"
$ws.Range("K27").Value = 0.3205128205128205
$ws.Range("L27").Value = 25
$ws.Range("M27").Value = 25
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 53

$ws.Range("J28").Value = "please"
$ws.Range("K28").Value = 0.3096234309623431
$ws.Range("L28").Value = 74
$ws.Range("M28").Value = 74
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 165
